$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 150.46153
$ws.Range("I55").Value = 182.2
$ws.Range("J55").Value = 130.625
$ws.Range("K55").Value = 182.2
$ws.Range("L55").Value = 130.625
$ws.Range("M55").Value = 31.80000000000001
$ws.Range("N55").Value = -558.625

# Row 69
$ws.Range("H69").Value = 3601
$ws.Range("I69").Value = 3156.5
$ws.Range("J69").Value = 3749.1667
$ws.Range("K69").Value = 9469.5
$ws.Range("L69").Value = 11247.5001
$ws.Range("M69").Value = -8595.5
$ws.Range("N69").Value = -12995.5001

# Row 70
$ws.Range("H70").Value = 1613.2667
$ws.Range("I70").Value = 1259.8
$ws.Range("J70").Value = 1790
$ws.Range("K70").Value = 3779.4
$ws.Range("L70").Value = 5370
$ws.Range("M70").Value = -3509.4
$ws.Range("N70").Value = -5910

# Row 72
$ws.Range("H72").Value = 3601
$ws.Range("I72").Value = 3156.5
$ws.Range("J72").Value = 3749.1667
$ws.Range("K72").Value = 28408.5
$ws.Range("L72").Value = 33742.5003
$ws.Range("M72").Value = -24040.5
$ws.Range("N72").Value = -42478.5003

# Row 73
$ws.Range("H73").Value = 1613.2667
$ws.Range("I73").Value = 1259.8
$ws.Range("J73").Value = 1790
$ws.Range("K73").Value = 3779.4
$ws.Range("L73").Value = 5370
$ws.Range("M73").Value = -2843.4
$ws.Range("N73").Value = -7242

# Row 74
$ws.Range("H74").Value = 3366.6667
$ws.Range("I74").Value = 2755.5557
$ws.Range("J74").Value = 3825
$ws.Range("K74").Value = 2755.5557
$ws.Range("L74").Value = 3825
$ws.Range("M74").Value = -1819.5557
$ws.Range("N74").Value = -5697

# Row 77
$ws.Range("H77").Value = 3366.6667
$ws.Range("I77").Value = 2755.5557
$ws.Range("J77").Value = 3825
$ws.Range("K77").Value = 13777.7785
$ws.Range("L77").Value = 19125
$ws.Range("M77").Value = -9097.7785
$ws.Range("N77").Value = -28485

# Row 100
$ws.Range("H100").Value = 1537.2106
$ws.Range("I100").Value = 1099.5
$ws.Range("J100").Value = 2023.5555
$ws.Range("K100").Value = 1099.5
$ws.Range("L100").Value = 2023.5555
$ws.Range("M100").Value = -558.5
$ws.Range("N100").Value = -3105.5555

# Row 105
$ws.Range("H105").Value = 39000
$ws.Range("J105").Value = 39000
$ws.Range("L105").Value = 39000
$ws.Range("N105").Value = -45988

# Row 132
$ws.Range("H132").Value = 2176.6924
$ws.Range("I132").Value = 1638.6389
$ws.Range("J132").Value = 8633.333000000001
$ws.Range("K132").Value = 4915.9167
$ws.Range("L132").Value = 25899.999
$ws.Range("M132").Value = -2385.9167
$ws.Range("N132").Value = -30959.999

# Row 138
$ws.Range("H138").Value = 2093.4395
$ws.Range("I138").Value = 1079
$ws.Range("J138").Value = 3868.7083
$ws.Range("K138").Value = 3237
$ws.Range("L138").Value = 11606.1249
$ws.Range("M138").Value = 1903
$ws.Range("N138").Value = -21886.1249

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7416.385
$ws.Range("I32").Value = 7144.7607
$ws.Range("J32").Value = 10171.429
$ws.Range("K32").Value = 7144.7607
$ws.Range("L32").Value = 10171.429
$ws.Range("M32").Value = -6857.7607
$ws.Range("N32").Value = -10745.429

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# Row 122
$ws.Range("H122").Value = 4273.6855
$ws.Range("I122").Value = 4377.9033
$ws.Range("K122").Value = 13133.7099
$ws.Range("M122").Value = -10683.7099

# Row 132
$ws.Range("H132").Value = 5815860
$ws.Range("I132").Value = 7814135.5
$ws.Range("J132").Value = 2694.3635
$ws.Range("K132").Value = 23442406.5
$ws.Range("L132").Value = 8083.0905
$ws.Range("M132").Value = -23439876.5
$ws.Range("N132").Value = -13143.0905

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1055
$ws.Range("I99").Value = 1061
$ws.Range("J99").Value = 1025
$ws.Range("K99").Value = 1061
$ws.Range("L99").Value = 1025
$ws.Range("M99").Value = 437
$ws.Range("N99").Value = -4021

# Row 134
$ws.Range("H134").Value = 2969.4
$ws.Range("I134").Value = 1816.75
$ws.Range("J134").Value = 7580
$ws.Range("K134").Value = 5450.25
$ws.Range("L134").Value = 22740
$ws.Range("M134").Value = -2915.25
$ws.Range("N134").Value = -27810

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 929.8570999999999
$ws.Range("I16").Value = 898.8261
$ws.Range("J16").Value = 1072.6
$ws.Range("K16").Value = 898.8261
$ws.Range("L16").Value = 1072.6
$ws.Range("M16").Value = -611.8261
$ws.Range("N16").Value = -1646.6

# Row 58
$ws.Range("H58").Value = 2233.3809
$ws.Range("I58").Value = 767.9231
$ws.Range("J58").Value = 4614.75
$ws.Range("K58").Value = 767.9231
$ws.Range("L58").Value = 4614.75
$ws.Range("M58").Value = -564.9231
$ws.Range("N58").Value = -5020.75

# Row 113
$ws.Range("H113").Value = 929.8570999999999
$ws.Range("I113").Value = 898.8261
$ws.Range("J113").Value = 1072.6
$ws.Range("K113").Value = 898.8261
$ws.Range("L113").Value = 1072.6
$ws.Range("M113").Value = 1271.1739
$ws.Range("N113").Value = -5412.6

# Row 134
$ws.Range("H134").Value = 1408.1951
$ws.Range("I134").Value = 1368.0667
$ws.Range("J134").Value = 1517.6364
$ws.Range("K134").Value = 4104.2001
$ws.Range("L134").Value = 4552.9092
$ws.Range("M134").Value = -1569.2001
$ws.Range("N134").Value = -9622.9092

# Row 136
$ws.Range("H136").Value = 2233.3809
$ws.Range("I136").Value = 767.9231
$ws.Range("J136").Value = 4614.75
$ws.Range("K136").Value = 2303.7693
$ws.Range("L136").Value = 13844.25
$ws.Range("M136").Value = 246.2307000000001
$ws.Range("N136").Value = -18944.25

$ws = $wb.Worksheets.Item("CUL")
# Row 42
$ws.Range("H42").Value = 5000
$ws.Range("J42").Value = 5000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -16068

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 10007
$ws.Range("J21").Value = 10007
$ws.Range("L21").Value = 10007
$ws.Range("N21").Value = -10353

# Row 30
$ws.Range("H30").Value = 10007
$ws.Range("J30").Value = 10007
$ws.Range("L30").Value = 10007
$ws.Range("N30").Value = -10217

# Row 102
$ws.Range("H102").Value = 3756.1667
$ws.Range("I102").Value = 3807.0476
$ws.Range("J102").Value = 3400
$ws.Range("K102").Value = 3807.0476
$ws.Range("L102").Value = 3400
$ws.Range("M102").Value = -2185.0476
$ws.Range("N102").Value = -6644

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 132
$ws.Range("H132").Value = 3621.1365
$ws.Range("I132").Value = 2684.9375
$ws.Range("J132").Value = 6117.6665
$ws.Range("K132").Value = 8054.8125
$ws.Range("L132").Value = 18352.9995
$ws.Range("M132").Value = -5524.8125
$ws.Range("N132").Value = -23412.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4923.969
$ws.Range("I7").Value = 4812
$ws.Range("J7").Value = 5216.3335
$ws.Range("K7").Value = 4812
$ws.Range("L7").Value = 5216.3335
$ws.Range("M7").Value = -4700
$ws.Range("N7").Value = -5440.3335

# Row 40
$ws.Range("H40").Value = 4390.0347
$ws.Range("I40").Value = 8714.286
$ws.Range("J40").Value = 3014.1365
$ws.Range("K40").Value = 8714.286
$ws.Range("L40").Value = 3014.1365
$ws.Range("M40").Value = -8578.286
$ws.Range("N40").Value = -3286.1365

# Row 46
$ws.Range("H46").Value = 1144.4445
$ws.Range("I46").Value = 1400
$ws.Range("J46").Value = 1071.4286
$ws.Range("K46").Value = 1400
$ws.Range("L46").Value = 1071.4286
$ws.Range("M46").Value = -1212
$ws.Range("N46").Value = -1447.4286

# Row 55
$ws.Range("H55").Value = 977.7778
$ws.Range("I55").Value = 400
$ws.Range("K55").Value = 400
$ws.Range("M55").Value = -227

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 126
$ws.Range("H126").Value = 4923.969
$ws.Range("I126").Value = 4812
$ws.Range("J126").Value = 5216.3335
$ws.Range("K126").Value = 14436
$ws.Range("L126").Value = 15649.0005
$ws.Range("M126").Value = -11966
$ws.Range("N126").Value = -20589.0005

# Row 132
$ws.Range("H132").Value = 5940.2036
$ws.Range("I132").Value = 3070
$ws.Range("J132").Value = 11133.904
$ws.Range("K132").Value = 9210
$ws.Range("L132").Value = 33401.712
$ws.Range("M132").Value = -6680
$ws.Range("N132").Value = -38461.712

# Row 136
$ws.Range("H136").Value = 7941599.5
$ws.Range("I136").Value = 11629692
$ws.Range("J136").Value = 12200.25
$ws.Range("K136").Value = 34889076
$ws.Range("L136").Value = 36600.75
$ws.Range("M136").Value = -34886526
$ws.Range("N136").Value = -41700.75

$ws = $wb.Worksheets.Item("WVR")
# Row 86
$ws.Range("H86").Value = 7825
$ws.Range("J86").Value = 7825
$ws.Range("L86").Value = 7825
$ws.Range("N86").Value = -10071

# Row 89
$ws.Range("H89").Value = 7825
$ws.Range("J89").Value = 7825
$ws.Range("L89").Value = 39125
$ws.Range("N89").Value = -50357

# Row 105
$ws.Range("H105").Value = 615
$ws.Range("J105").Value = 615
$ws.Range("L105").Value = 615
$ws.Range("N105").Value = -7603

# Row 111
$ws.Range("H111").Value = 39000
$ws.Range("J111").Value = 39000
$ws.Range("L111").Value = 39000
$ws.Range("N111").Value = -47180

# Row 112
$ws.Range("H112").Value = 18225
$ws.Range("J112").Value = 18225
$ws.Range("L112").Value = 18225
$ws.Range("N112").Value = -21179

# Row 113
$ws.Range("H113").Value = 1182.8334
$ws.Range("I113").Value = 420.81818
$ws.Range("J113").Value = 2380.2856
$ws.Range("K113").Value = 1262.45454
$ws.Range("L113").Value = 7140.8568
$ws.Range("M113").Value = 907.54546
$ws.Range("N113").Value = -11480.8568

# Row 122
$ws.Range("H122").Value = 1618.0769
$ws.Range("I122").Value = 1890.7894
$ws.Range("J122").Value = 877.8570999999999
$ws.Range("K122").Value = 5672.3682
$ws.Range("L122").Value = 2633.5713
$ws.Range("M122").Value = -3222.3682
$ws.Range("N122").Value = -7533.5713

# Row 132
$ws.Range("H132").Value = 1398.4546
$ws.Range("I132").Value = 1168.2075
$ws.Range("J132").Value = 2337.1538
$ws.Range("K132").Value = 3504.6225
$ws.Range("L132").Value = 7011.4614
$ws.Range("M132").Value = -974.6224999999999
$ws.Range("N132").Value = -12071.4614

# Row 136
$ws.Range("H136").Value = 601.9474
$ws.Range("I136").Value = 560.6458
$ws.Range("J136").Value = 822.2222
$ws.Range("K136").Value = 1681.9374
$ws.Range("L136").Value = 2466.6666
$ws.Range("M136").Value = 868.0626
$ws.Range("N136").Value = -7566.6666

Write-Output "Applied Ultima_Profits updates"